$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14
$ws.Cells.Item($row, 1).Value = "2025-08-15 04:03:17 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-15 09:33:17 IST"
$ws.Cells.Item($row, 3).Value = "UPDATED"
$ws.Cells.Item($row, 4).Value = "New circular processed."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = "INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 2

$ws.Range("A13:H13").Copy()
$ws.Range("A14:H14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
